$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5241678953170776
$ws.Range("B1").Value = 3.208840370178223
$ws.Range("C1").Value = 5.945469856262207
$ws.Range("D1").Value = 1.48795473575592
$ws.Range("E1").Value = 0.8697299957275391
